# Nova Priority Guide: split the "playable but annoying..." run so that
# "may be " is inserted between "playable but " and "annoying, ...".
#
# A plain text replace in this runtime re-merges every run in the
# paragraph (it collapses same-format adjacent runs back into one run
# whenever any text in the paragraph changes), which would lose the
# three-way run split the target OOXML expects. Turning on
# TrackRevisions around the insertion keeps the new text in its own run
# (wrapped in <w:ins>) instead of being folded back in; accepting just
# that revision afterwards then bakes the insertion in as a clean,
# separate <w:r> run with no leftover formatting markup, and without
# touching unrelated runs elsewhere in the document the way accepting
# every revision in the document would.

$d = $word.ActiveDocument

$needle = "playable but annoying, acceptable for a Beta version"
$splitAfter = "playable but "
$insertion = "may be "

$fullText = $d.Content.Text
$needleStart = $fullText.IndexOf($needle)
if ($needleStart -lt 0) {
    throw "Could not find target text '$needle' in document"
}

$insertAt = $needleStart + $splitAfter.Length
$insertionPoint = $d.Range($insertAt, $insertAt)

$wasTrackingRevisions = $d.TrackRevisions

$d.TrackRevisions = $true
$insertionPoint.InsertBefore($insertion)
$d.TrackRevisions = $wasTrackingRevisions

$d.Revisions.Item(1).Accept()
